# Update EventRanking values on the "Events" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# Column C holds "EventRanking"; update the ranking values for rows 2-7
$ws.Range("C2").Value = 4
$ws.Range("C3").Value = 7
$ws.Range("C4").Value = 6
$ws.Range("C5").Value = 10
$ws.Range("C6").Value = 5
$ws.Range("C7").Value = 1

# Move the active selection to C5, matching the saved view state
$ws.Activate()
$ws.Range("C5").Select()
